# [DOC_02] Change the navigation tag for news items from "News" to "HR" (RH),
# since a news item cannot be tagged with the generic "News" category but
# must use a sub-category instead.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "HR"

# Reflect the new active selection captured in the saved workbook (E2).
$ws.Range("E2").Select()
